$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column L (year 2020), mirroring column K (year 2019) row by row
$ws.Range("L4").Value = 2020
$ws.Range("L5").Value = 5.6
$ws.Range("L6").Value = 0.8
$ws.Range("L7").Value = 1.9
$ws.Range("L8").Value = 0.7
$ws.Range("L9").Value = 0.7
$ws.Range("L10").Value = 0.9
$ws.Range("L11").Value = 0.3
$ws.Range("L12").Value = 0.2

# Copy the number formatting/style from column K over to the new column L
$ws.Range("K4:K12").Copy()
$ws.Range("L4:L12").PasteSpecial(-4122) # xlPasteFormats

# Update the selection to match the target state
$ws.Range("N5").Select()
